$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11:C37").ClearContents()
$ws.Range("A11:C37").Font.Name = "Arial"
